$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 19097
$ws.Range("J19").Value = 20806.7
$ws.Range("L19").Value = 20806.7
$ws.Range("N19").Value = -21156.7

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1809.7273
$ws.Range("J43").Value = 1711.7778
$ws.Range("L43").Value = 1711.7778
$ws.Range("N43").Value = -1849.7778

# Sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 4007.25
$ws.Range("I111").Value = 10029
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 30087
$ws.Range("L111").Value = 6000
$ws.Range("M111").Value = -27020
$ws.Range("N111").Value = -12134

# Sheet ALC, row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1301.3513
$ws.Range("I135").Value = 487.5484
$ws.Range("K135").Value = 4387.9356
$ws.Range("M135").Value = -1852.9356

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1793.7073
$ws.Range("I137").Value = 1775.9487
$ws.Range("J137").Value = 2140
$ws.Range("K137").Value = 5327.8461
$ws.Range("L137").Value = 6420
$ws.Range("M137").Value = -2777.8461
$ws.Range("N137").Value = -11520

# Sheet ARM, row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 11827.579
$ws.Range("J23").Value = 8595.444
$ws.Range("L23").Value = 8595.444
$ws.Range("N23").Value = -9113.444

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 889040.2
$ws.Range("I32").Value = 1033978.7
$ws.Range("K32").Value = 1033978.7
$ws.Range("M32").Value = -1033691.7

# Sheet ARM, row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 12021.286
$ws.Range("J37").Value = 12021.286
$ws.Range("L37").Value = 12021.286
$ws.Range("N37").Value = -12567.286

# Sheet ARM, row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 29408.8
$ws.Range("J44").Value = 29408.8
$ws.Range("L44").Value = 29408.8
$ws.Range("N44").Value = -30384.8

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4603.5
$ws.Range("I45").Value = 5900
$ws.Range("J45").Value = 3307
$ws.Range("K45").Value = 5900
$ws.Range("L45").Value = 3307
$ws.Range("M45").Value = -5523
$ws.Range("N45").Value = -4061

# Sheet ARM, row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 35520.137
$ws.Range("I63").Value = 183319.8
$ws.Range("K63").Value = 183319.8
$ws.Range("M63").Value = -182633.8

# Sheet ARM, row 64
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Sheet ARM, row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 35520.137
$ws.Range("I66").Value = 183319.8
$ws.Range("K66").Value = 916599
$ws.Range("M66").Value = -913167

# Sheet ARM, row 67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2464.2983
$ws.Range("I132").Value = 1390.7333
$ws.Range("K132").Value = 4172.199900000001
$ws.Range("M132").Value = -1642.199900000001

# Sheet BSM, row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 277.14285
$ws.Range("I22").Value = 277.14285
$ws.Range("K22").Value = 277.14285
$ws.Range("M22").Value = -104.14285

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2016.1968
$ws.Range("I134").Value = 1582.4
$ws.Range("J134").Value = 3236.25
$ws.Range("K134").Value = 4747.200000000001
$ws.Range("L134").Value = 9708.75
$ws.Range("M134").Value = -2212.200000000001
$ws.Range("N134").Value = -14778.75

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4215.8843
$ws.Range("I31").Value = 978.0465
$ws.Range("J31").Value = 9570.77
$ws.Range("K31").Value = 978.0465
$ws.Range("L31").Value = 9570.77
$ws.Range("M31").Value = -683.0465
$ws.Range("N31").Value = -10160.77

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4215.8843
$ws.Range("I34").Value = 978.0465
$ws.Range("J34").Value = 9570.77
$ws.Range("K34").Value = 978.0465
$ws.Range("L34").Value = 9570.77
$ws.Range("M34").Value = -776.0465
$ws.Range("N34").Value = -9974.77

# Sheet CRP, row 56
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 3093
$ws.Range("I56").Value = 3093
$ws.Range("K56").Value = 3093
$ws.Range("M56").Value = -2248

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1064.9791
$ws.Range("I58").Value = 811.5161000000001
$ws.Range("J58").Value = 1527.1765
$ws.Range("K58").Value = 811.5161000000001
$ws.Range("L58").Value = 1527.1765
$ws.Range("M58").Value = -608.5161000000001
$ws.Range("N58").Value = -1933.1765

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1064.9791
$ws.Range("I136").Value = 811.5161000000001
$ws.Range("J136").Value = 1527.1765
$ws.Range("K136").Value = 2434.5483
$ws.Range("L136").Value = 4581.529500000001
$ws.Range("M136").Value = 115.4516999999996
$ws.Range("N136").Value = -9681.529500000001

# Sheet CUL, row 93
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 4057.3845
$ws.Range("I93").Value = 900
$ws.Range("J93").Value = 4320.5
$ws.Range("K93").Value = 2700
$ws.Range("L93").Value = 12961.5
$ws.Range("M93").Value = -828
$ws.Range("N93").Value = -16705.5

# Sheet CUL, row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7943661.5
$ws.Range("J137").Value = 4350
$ws.Range("L137").Value = 13050
$ws.Range("N137").Value = -23250

# Sheet GSM, row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 45.6
$ws.Range("I2").Value = 45.6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 45.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 67.40000000000001
$ws.Range("N2").ClearContents()

# Sheet GSM, row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576

# Sheet GSM, row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 10000
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -11996

# Sheet GSM, row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 10000
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -39984

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1995.6086
$ws.Range("I46").Value = 1989.421
$ws.Range("J46").Value = 2025
$ws.Range("K46").Value = 1989.421
$ws.Range("L46").Value = 2025
$ws.Range("M46").Value = -1801.421
$ws.Range("N46").Value = -2401

# Sheet LTW, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 426.625
$ws.Range("I55").Value = 329.75
$ws.Range("J55").Value = 523.5
$ws.Range("K55").Value = 329.75
$ws.Range("L55").Value = 523.5
$ws.Range("M55").Value = -156.75
$ws.Range("N55").Value = -869.5

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4631289.5
$ws.Range("I136").Value = 1687.3667
$ws.Range("K136").Value = 5062.1001
$ws.Range("M136").Value = -2512.1001

# Sheet WVR, row 86
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 62500
$ws.Range("J86").Value = 62500
$ws.Range("L86").Value = 62500
$ws.Range("N86").Value = -64746

# Sheet WVR, row 89
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value = 62500
$ws.Range("J89").Value = 62500
$ws.Range("L89").Value = 312500
$ws.Range("N89").Value = -323732

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1138.1111
$ws.Range("I107").Value = 1092.875
$ws.Range("K107").Value = 3278.625
$ws.Range("M107").Value = -1358.625

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2454.6667
$ws.Range("I113").Value = 3100.3333
$ws.Range("J113").Value = 1163.3334
$ws.Range("K113").Value = 9300.999899999999
$ws.Range("L113").Value = 3490.0002
$ws.Range("M113").Value = -7130.999899999999
$ws.Range("N113").Value = -7830.0002

# Sheet WVR, row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 31785.8
$ws.Range("J123").Value = 49464.5
$ws.Range("L123").Value = 49464.5
$ws.Range("N123").Value = -59264.5

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1397.9508
$ws.Range("I132").Value = 1251.3334
$ws.Range("K132").Value = 3754.0002
$ws.Range("M132").Value = -1224.0002
